# Update "想去人数" (want-to-go count, column F) for a handful of events
# on the "展览" and "全部类型" sheets, per the upstream data refresh
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F on sheet "展览"
$exhibitionUpdates = @{
    4  = 1153
    5  = 132
    10 = 16854
    12 = 215
    14 = 6495
    20 = 1284
    21 = 96
    23 = 646
    29 = 78
    31 = 517
    33 = 11536
    34 = 1260
    38 = 3869
}

# Row -> new value for column F on sheet "全部类型"
$allTypesUpdates = @{
    4  = 1153
    5  = 132
    10 = 16854
    12 = 215
    14 = 6495
    20 = 1284
    21 = 96
    23 = 646
    29 = 78
    31 = 517
    34 = 11536
    35 = 1260
    39 = 3869
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
